# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# per the latest scrape. D-column writes are wrapped with a
# temporary Text number-format so numeric-looking price strings
# (e.g. "566.87") are stored as text, matching the source data,
# then the cell style is reset to Normal to avoid stray formatting diffs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.209.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.625.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +4.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.649.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  +6.16%  "
$ws.Range("E12").Value = "  +7.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.094.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.225.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.73%  "
$ws.Range("E17").Value = "  +5.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.643.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.438"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.24%  "
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0802"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.12%  "
$ws.Range("E31").Value = "  +4.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.887"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.62%  "
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.885"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "299.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0984"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0544"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("E51").Value = "  +7.58%  "
